$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.844.54'
$ws.Range("E2").Value = '  +0.12%  '
$ws.Range("D3").Value = '3.503.57'
$ws.Range("E3").Value = '  -1.05%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '606.47'
$ws.Range("E5").Value = '  -0.82%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '150.23'
$ws.Range("E6").Value = '  -1.41%  '
$ws.Range("D7").Value = '3.499.90'
$ws.Range("E7").Value = '  -1.13%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  +0.76%  '
$ws.Range("E10").Value = '  +3.00%  '
$ws.Range("E11").Value = '  +6.94%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.430'
$ws.Range("E12").Value = '  +0.90%  '
$ws.Range("E13").Value = '  -1.72%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '31.96'
$ws.Range("E14").Value = '  +0.07%  '
$ws.Range("D15").Value = '4.095.49'
$ws.Range("E15").Value = '  -1.10%  '
$ws.Range("D16").Value = '67.766.47'
$ws.Range("E16").Value = '  +0.29%  '
$ws.Range("D17").Value = '3.501.69'
$ws.Range("E17").Value = '  -1.59%  '
$ws.Range("E18").Value = '  -0.18%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.47'
$ws.Range("E19").Value = '  +1.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.36'
$ws.Range("E20").Value = '  +0.92%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.93'
$ws.Range("E21").Value = '  +2.51%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '444.84'
$ws.Range("E22").Value = '  -0.29%  '
$ws.Range("E23").Value = '  +0.26%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.96'
$ws.Range("E24").Value = '  +2.61%  '
$ws.Range("D25").Value = '3.644.31'
$ws.Range("E25").Value = '  -1.06%  '
$ws.Range("E26").Value = '  -0.10%  '
$ws.Range("E27").Value = '  -3.15%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.68'
$ws.Range("E28").Value = '  -3.05%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.99'
$ws.Range("E29").Value = '  -1.88%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.67'
$ws.Range("E30").Value = '  +0.76%  '
$ws.Range("E31").Value = '  -1.33%  '
$ws.Range("E32").Value = '  +1.77%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.170'
$ws.Range("E33").Value = '  +1.16%  '
$ws.Range("E34").Value = '  -0.20%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '25.57'
$ws.Range("E35").Value = '  -0.31%  '
$ws.Range("E36").Value = '  -0.49%  '
$ws.Range("E37").Value = '  +0.16%  '
$ws.Range("D38").Value = '3.499.49'
$ws.Range("E38").Value = '  -0.81%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.97'
$ws.Range("E39").Value = '  -0.80%  '
$ws.Range("E40").Value = '  -0.02%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.33'
$ws.Range("E41").Value = '  +6.40%  '
$ws.Range("E42").Value = '  -0.06%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '176.49'
$ws.Range("E43").Value = '  +0.06%  '
$ws.Range("E44").Value = '  +0.64%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.41'
$ws.Range("E45").Value = '  +0.36%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.896'
$ws.Range("E46").Value = '  +0.99%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '30.19'
$ws.Range("E47").Value = '  +5.07%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '46.75'
$ws.Range("E48").Value = '  +2.63%  '
$ws.Range("E49").Value = '  +1.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.52'
$ws.Range("E50").Value = '  -5.36%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.61'
$ws.Range("E51").Value = '  +0.13%  '
